# Update "想去人数" (column F) values on sheets "展览" and "全部类型"
# to match newly scraped counts.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 165
    4  = 75
    6  = 536
    7  = 1632
    9  = 21
    10 = 1460
    12 = 41
    13 = 370
    15 = 186
    17 = 17
    19 = 261
    22 = 201
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}

$updates2 = @{
    3  = 165
    4  = 75
    6  = 536
    7  = 1632
    10 = 21
    11 = 1460
    13 = 41
    14 = 370
    16 = 186
    18 = 17
    20 = 261
    23 = 201
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates2.Keys) {
    $ws4.Range("F$row").Value = $updates2[$row]
}
